$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (sending cluster, ligand symbol, receptor symbol, target cluster, then 16 numeric metrics)
# Clusters now include "ECs" in addition to the previously existing "FAPs" and "sCs".
$rows = @(
    @("ECs",  "Clcf1", "Cntfr", "ECs",  2, 0.6666666666666666, 1.798783666666667, 5.396351,          0.2319744053785674, 0.2319744053785674, 1, 0.3333333333333333, 0.02507166666666667, 0.075215,   0.001520672841173258, 0.001520672841173258, 0.04509850449611112, 0.4058865404650001, 0.0003527571781065032, 0.0003527571781065032),
    @("ECs",  "Clcf1", "Cntfr", "FAPs", 2, 0.6666666666666666, 1.798783666666667, 5.396351,          0.2319744053785674, 0.2319744053785674, 3, 1,                   13.42247233333333,   40.267417,  0.8141137727328104,   0.8141137727328106,   24.14412399948522,   217.297115995367,  0.1888535583401959,   0.1888535583401959),
    @("ECs",  "Clcf1", "Cntfr", "sCs",  2, 0.6666666666666666, 1.798783666666667, 5.396351,          0.2319744053785674, 0.2319744053785674, 3, 1,                   3.039675333333333,   9.119026,   0.1843655544260162,   0.1843655544260162,   5.467718341569555,   49.209465074126,   0.04276808986026501,  0.04276808986026501),
    @("FAPs", "Clcf1", "Cntfr", "ECs",  3, 1,                   1.782700333333333, 5.348101,          0.2299002695301921, 0.2299002695301921, 1, 0.3333333333333333, 0.02507166666666667, 0.075215,   0.001520672841173258, 0.001520672841173258, 0.04469526852388889, 0.402257416715,    0.000349603096052975, 0.000349603096052975),
    @("FAPs", "Clcf1", "Cntfr", "FAPs", 3, 1,                   1.782700333333333, 5.348101,          0.2299002695301921, 0.2299002695301921, 3, 1,                   13.42247233333333,   40.267417,  0.8141137727328104,   0.8141137727328106,   23.92824590279078,   215.354213125117,  0.1871649757795147,   0.1871649757795147),
    @("FAPs", "Clcf1", "Cntfr", "sCs",  3, 1,                   1.782700333333333, 5.348101,          0.2299002695301921, 0.2299002695301921, 3, 1,                   3.039675333333333,   9.119026,   0.1843655544260162,   0.1843655544260162,   5.418830229958444,   48.769472069626,   0.04238569065462442,  0.04238569065462442),
    @("sCs",  "Clcf1", "Cntfr", "ECs",  3, 1,                   4.172749333333334, 12.518248,         0.5381253250912406, 0.5381253250912404, 1, 0.3333333333333333, 0.02507166666666667, 0.075215,   0.001520672841173258, 0.001520672841173258, 0.1046177803688889,  0.9415600233200002, 0.0008183125670137798, 0.0008183125670137798),
    @("sCs",  "Clcf1", "Cntfr", "FAPs", 3, 1,                   4.172749333333334, 12.518248,         0.5381253250912406, 0.5381253250912404, 3, 1,                   13.42247233333333,   40.267417,  0.8141137727328104,   0.8141137727328106,   56.00861248060178,   504.0775123254161, 0.4380952386130999,   0.4380952386131),
    @("sCs",  "Clcf1", "Cntfr", "sCs",  3, 1,                   4.172749333333334, 12.518248,         0.5381253250912406, 0.5381253250912404, 3, 1,                   3.039675333333333,   9.119026,   0.1843655544260162,   0.1843655544260162,   12.68380322071645,   114.154228986448,  0.09921177391112676,  0.09921177391112676)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $rowData[$c - 1]
    }
}
